$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "5.01") are stored as text, matching the original inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.065.84"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.338.41"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D5").Value = "303.48"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "95.15"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").Value = "34.24"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").Value = "19.06"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "0.0786"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").Value = "6.72"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").Value = "2.693.67"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "2.333.03"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "42.986.32"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "12.12"
$ws.Range("E19").Value = "  -3.73%  "
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "0.0₃0893"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "68.08"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "236.90"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "2.25"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "24.72"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").Value = "  -7.29%  "
$ws.Range("D29").Value = "9.14"
$ws.Range("D30").Value = "31.67"
$ws.Range("E30").Value = "  -4.35%  "
$ws.Range("D31").Value = "141.70"
$ws.Range("E31").Value = "  -14.63%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "5.01"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("D35").Value = "0.0702"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").Value = "4.41"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "22.39"
$ws.Range("E40").Value = "  +23.15%  "
$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").Value = "1.941.39"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("D45").Value = "10.10"
$ws.Range("E45").Value = "  -5.56%  "
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "2.74"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").Value = "2.560.05"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "52.81"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").Value = "72.63"
$ws.Range("E51").Value = "  +0.78%  "

# Restore default (Normal) style so no stray number-format style sticks
# around on cells, matching the original workbook styling.
$textRange.Style = "Normal"
